$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full data set for rows 2-10 (3x3 combinations of sending/target clusters ECs/FAPs/sCs)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics
$data = @(
    ,("ECs", "Reln", "Itgb1", "ECs", 2, 0.6666666666666666, 0.2638773333333334, 0.791632, 0.0472992197617186, 0.04729921976171859, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 29.689733845248, 267.207604607232, 0.01549197483049645, 0.01549197483049645)
    ,("ECs", "Reln", "Itgb1", "FAPs", 2, 0.6666666666666666, 0.2638773333333334, 0.791632, 0.0472992197617186, 0.04729921976171859, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 28.05397778283734, 252.485800045536, 0.01463844438526634, 0.01463844438526634)
    ,("ECs", "Reln", "Itgb1", "sCs", 2, 0.6666666666666666, 0.2638773333333334, 0.791632, 0.0472992197617186, 0.04729921976171859, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 32.90330149827911, 296.129713484512, 0.0171688005459558, 0.0171688005459558)
    ,("FAPs", "Reln", "Itgb1", "ECs", 3, 1, 1.127914, 3.383742, 0.2021751981665183, 0.2021751981665183, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 126.905430024288, 1142.148870218592, 0.06621870502568583, 0.06621870502568583)
    ,("FAPs", "Reln", "Itgb1", "FAPs", 3, 1, 1.127914, 3.383742, 0.2021751981665183, 0.2021751981665183, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 119.913574603924, 1079.222171435316, 0.06257038507929176, 0.06257038507929177)
    ,("FAPs", "Reln", "Itgb1", "sCs", 3, 1, 1.127914, 3.383742, 0.2021751981665183, 0.2021751981665183, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 140.6414637336413, 1265.773173602772, 0.07338610806154067, 0.07338610806154068)
    ,("sCs", "Reln", "Itgb1", "ECs", 3, 1, 4.187102666666667, 12.561308, 0.7505255820717631, 0.750525582071763, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 471.105123678912, 4239.946113110208, 0.2458206178806741, 0.245820617880674)
    ,("sCs", "Reln", "Itgb1", "FAPs", 3, 1, 4.187102666666667, 12.561308, 0.7505255820717631, 0.750525582071763, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 445.1495840938427, 4006.346256844584, 0.2322771294796082, 0.2322771294796082)
    ,("sCs", "Reln", "Itgb1", "sCs", 3, 1, 4.187102666666667, 12.561308, 0.7505255820717631, 0.750525582071763, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 522.096762557281, 4698.870863015528, 0.2724278347114808, 0.2724278347114808)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $data) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $row[$i]
    }
    $r = $r + 1
}
